$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add tested count for row 35 (Apr 9 data correction)
$ws.Range("B35").Value = 920

# Add new row 36 for Apr 10 data
$ws.Range("A36").Value = 43930
$ws.Range("C36").Value = 1628
$ws.Range("D36").Value = 10
$ws.Range("F36").Value = 20
$ws.Range("H36").Value = 255
$ws.Range("J36").Value = 320
$ws.Range("L36").Value = 283
$ws.Range("N36").Value = 293
$ws.Range("P36").Value = 224
$ws.Range("R36").Value = 131
$ws.Range("T36").Value = 90
$ws.Range("V36").Value = 2
$ws.Range("X36").Value = 790
$ws.Range("Y36").Value = 832
$ws.Range("Z36").Value = 6
$ws.Range("AA36").Value = 348
$ws.Range("AB36").Value = 132
$ws.Range("AC36").Value = 40
$ws.Range("AD36").Value = 43
$ws.Range("AE36").Value = 154
$ws.Range("AF36").Value = 5
$ws.Range("AG36").Value = 8
$ws.Range("AH36").Value = 93
$ws.Range("AI36").Value = 32
$ws.Range("AJ36").Value = 39
$ws.Range("AK36").Value = 10
$ws.Range("AL36").Value = 30
$ws.Range("AM36").Value = 14
$ws.Range("AN36").Value = 32
$ws.Range("AO36").Value = 34
$ws.Range("AP36").Value = 16
$ws.Range("AQ36").Value = 821
$ws.Range("AR36").Value = 21
$ws.Range("AS36").Value = 17
$ws.Range("AT36").Value = 5
$ws.Range("AU36").Value = 26
$ws.Range("AV36").Value = 1
$ws.Range("AW36").Value = 13
$ws.Range("AX36").Value = 1
$ws.Range("AY36").Value = 1
$ws.Range("AZ36").Value = 7
$ws.Range("BA36").Value = 4
$ws.Range("BB36").Value = 13
$ws.Range("BC36").Value = 2
$ws.Range("BD36").Value = 10
$ws.Range("BE36").Value = 14
$ws.Range("BG36").Value = 39
$ws.Range("BH36").Value = 3
$ws.Range("BI36").Value = 97

# Reset view: scroll to top-left and select F17
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$null = $ws.Range("F17").Select()

Write-Output "Done"